$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (B1) into C1 before writing its value, so C1 matches
# the existing header formatting (bold, border, centered) exactly.
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)

$values = @{
    1 = 2
    2 = -0.3112980140467429
    3 = 0.2065416847512486
    4 = 0.1418205654855191
    5 = -0.07704090021470952
    6 = 0.1425763902569039
    7 = -0.4527558810974688
    8 = -0.2603289572191466
    9 = -0.4114738013147691
    10 = 0.3580686953613876
    11 = -0.2399055122828032
    12 = -0.0873066370082155
    13 = -0.01191385433495715
    14 = 0.05155802922678206
    15 = -0.01369572637148195
    16 = 0.3472103327010492
    17 = 0.598514994193861
    18 = 0.03878916555125164
    19 = 0.3958043559957307
    20 = 0.2934706472736356
    21 = 0.517878593557578
    22 = 0.377430769596674
    23 = -0.0608179811506077
    24 = 4.579137554806769
    25 = 0.488382454367825
    26 = 0.3615177383265028
    27 = 0.303409187322828
    28 = 1.069211801959468
    29 = 5.465989120357457
    30 = 0.9723697196820085
    31 = -0.1983294097244325
    32 = 0.7801612285857227
    33 = 0.912811341700061
    34 = -0.6563600928369668
    35 = 0.8136177197400545
    36 = 0.7676957230286918
    37 = 0.7417196672734563
    38 = 0.7176524242778116
    39 = 0.580898655261708
    40 = 0.7521349625251679
    41 = 0.5543715805352603
    42 = 0.6904290198476798
    43 = 0.7091434782750188
    44 = 0.674411792363956
    45 = 0.6580905976512383
    46 = -1.263811510086398
    47 = -0.9795097922288276
    48 = -0.8735394605770974
    49 = -0.6425263145464726
    50 = -0.0488315234424316
    51 = -0.8654141101463644
    52 = -0.8654141101463644
    53 = -1.097302253505303
    54 = -0.1871342512186376
    55 = -0.9993792754448271
    56 = -0.8844730449215408
    57 = -0.9341162932906286
    58 = -1.128542268251601
    59 = -0.8483927370979433
    60 = -0.4981670011442724
    61 = 0.3713187932361238
    62 = -1.209750680733107
    63 = -0.7400286323080764
    64 = -0.9027772512895821
    65 = -0.1119708445767711
    66 = -0.7911518241820487
    67 = -0.7628194820164587
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 3).Value = $values[$row]
}
